# Update the data row (row 2) on the active worksheet with new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 looks numeric ("99004392.0") but must stay stored as text, so force
# the cell to a text format before assigning the value.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "99004392.0"

$ws.Range("B2").Value = "midnight children"
$ws.Range("C2").Value = "catcher in rye"
$ws.Range("D2").Value = "animal farm"
